# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Papaya" (Vega Modelo de Temuco)
# above the existing row 108, shifting the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(108).Insert()

$ws.Cells.Item(108, 1).Value  = 10
$ws.Cells.Item(108, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(108, 3).Value  = "La Araucanía"
$ws.Cells.Item(108, 4).Value  = 45180
$ws.Cells.Item(108, 5).Value  = 9
$ws.Cells.Item(108, 6).Value  = "Fruta"
$ws.Cells.Item(108, 7).Value  = 100108
$ws.Cells.Item(108, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(108, 9).Value  = 100108004
$ws.Cells.Item(108, 10).Value = "Papaya"
$ws.Cells.Item(108, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(108, 12).Value = "Primera"
$ws.Cells.Item(108, 13).Value = 300
$ws.Cells.Item(108, 14).Value = 24000
$ws.Cells.Item(108, 15).Value = 24000
$ws.Cells.Item(108, 16).Value = 24000
$ws.Cells.Item(108, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(108, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(108, 19).Value = 2400
$ws.Cells.Item(108, 20).Value = 10
